# Applies the "added o to output of basic model" edit to Statistics.xlsx
# - Adds a new COIN-BC solver results block in columns F:H (values in column F)
# - Fills in the previously-missing K3/K4 values for the Model 2 / Gecode block
# - Replaces "DNF" text entries in column B with numeric results (-1 or 1.29)
# - Moves the active cell selection to K10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New COIN-BC header label (row 1) for the new results block in column F
$ws.Range("F1").Value = "COIN-BC"

# New COIN-BC timing data in column F
$ws.Range("F3").Value = 1460
$ws.Range("F4").Value = 16949
$ws.Range("F5").Value = 23878
$ws.Range("F6").Value = 809
$ws.Range("F7").Value = 32214
$ws.Range("F8").Value = 621
$ws.Range("F10").Value = 11844
$ws.Range("F11").Value = 4705
$ws.Range("F12").Value = 859
$ws.Range("F16").Value = 15934
$ws.Range("F18").Value = 3872

# Newly filled in data for the Model 2 / Gecode block
$ws.Range("K3").Value = 480
$ws.Range("K4").Value = 1230

# Column B (Model 1 / Gecode) results that were previously blank become -1
$ws.Range("B9").Value = -1
$ws.Range("B13").Value = -1
$ws.Range("B14").Value = -1
$ws.Range("B15").Value = -1
$ws.Range("B17").Value = -1
$ws.Range("B19").Value = -1
$ws.Range("B20").Value = -1
$ws.Range("B21").Value = -1

# Column B entries that used to hold the text "DNF" become numeric results
$ws.Range("B11").Value = 1.29
$ws.Range("B22").Value = -1

# Move the active selection
$ws.Range("K10").Select()
